$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert 7 new blank rows before the old "ItemDataTable" row (row 8), shifting
# everything below down by 7 (old row8 -> row15, old row10 -> row17, etc.)
$ws.Range("A8:A14").EntireRow.Insert()

# Row 8: new "UseDB" toggle
$ws.Range("A8").Value = "UseDB"
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = "Use auxiliary database for error codes and item data"
$ws.Range("C8").WrapText = $true

# Row 11 / Row 12: new UseDataServiceErrors / UseDataServiceItem toggles
$ws.Range("A11").Value = "UseDataServiceErrors"
$ws.Range("A12").Value = "UseDataServiceItem"

$ws.Range("C12").Value = "Use Data service for item data"
$ws.Range("C12").WrapText = $true

$ws.Range("C11").Value = "Use Data service for error codes"
$ws.Range("C11").WrapText = $true

$ws.Range("B11").Value = $true
$ws.Range("B12").Value = $true

# Row 14 / Row 13: ExceptionLanguage and EnforceSecondDataSource
$ws.Range("A14").Value = "ExceptionLanguage"
$ws.Range("A13").Value = "EnforceSecondDataSource"

$ws.Range("C13").Value = "If set to true, framework will attempt to get data from relevant data source (DB/Data Service). If data doesn't exist, an exception will be htrown . If false, it will not try to retrieve data from data source "
$ws.Range("C13").WrapText = $true

$ws.Range("B13").Value = $true
$ws.Range("B14").Value = "es-ES"

# Rows 15 / 16: blank description-styled placeholder rows
$ws.Range("C14").WrapText = $true
$ws.Range("C15").WrapText = $true
$ws.Range("C16").WrapText = $true

# Row 7: blank row with formatting only
$ws.Range("C7").WrapText = $true

# Update the active selection to match the authored state
$ws.Range("B10").Select()
